# Updates cryptos list values (Coin, Link, Price, Volume(1h)) for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, [string]$text) {
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

$data = @(
    @{ row=2; b='Bitcoin'; c='https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; d='64.208.47'; e='  +0.36%  ' }
    @{ row=3; b='Ethereum'; c='https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; d='3.155.59'; e='  -0.91%  ' }
    @{ row=4; b='TetherUSD'; c='https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; d='1.00'; e='  -0.05%  ' }
    @{ row=5; b='BNB'; c='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; d='567.36'; e='  -0.59%  ' }
    @{ row=6; b='Solana'; c='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; d='162.99'; e='  -4.05%  ' }
    @{ row=7; b='USDC'; c='https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; d='1.00'; e='  +0.01%  ' }
    @{ row=8; b='XRP'; c='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; d='0.584'; e='  -4.75%  ' }
    @{ row=9; b='Dogecoin'; c='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; d='0.117'; e='  -3.59%  ' }
    @{ row=10; b='Toncoin'; c='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; d='6.66'; e='  -1.37%  ' }
    @{ row=11; b='Cardano'; c='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; d='0.384'; e='  -0.48%  ' }
    @{ row=12; b='WrappedliquidstakedEther2.0'; c='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; d='3.719.46'; e='  -1.00%  ' }
    @{ row=13; b='TRON'; c='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; d='0.128'; e='  -0.92%  ' }
    @{ row=14; b='WrappedBTC'; c='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; d='64.413.07'; e='  +0.49%  ' }
    @{ row=15; b='Avalanche'; c='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; d='25.28'; e='  -0.67%  ' }
    @{ row=16; b='WrappedEther'; c='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; d='3.174.93'; e='  -0.53%  ' }
    @{ row=17; b='ShibaInu'; c='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; d='0.0000155'; e='  -2.00%  ' }
    @{ row=18; b='BitcoinCash'; c='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; d='405.27'; e='  -2.26%  ' }
    @{ row=19; b='Chainlink'; c='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; d='12.67'; e='  -1.03%  ' }
    @{ row=20; b='Polkadot'; c='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; d='5.26'; e='  -2.18%  ' }
    @{ row=21; b='Uniswap'; c='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; d='7.08'; e='  -1.19%  ' }
    @{ row=22; b='Dai'; c='https://coinranking.com/coin/MoTuySvg7+dai-dai'; d='1.00'; e='  +0.20%  ' }
    @{ row=23; b='Litecoin'; c='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; d='68.53'; e='  -3.60%  ' }
    @{ row=24; b='Kaspa'; c='https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; d='0.200'; e='  -1.50%  ' }
    @{ row=25; b='Polygon'; c='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; d='0.484'; e='  -1.72%  ' }
    @{ row=26; b='PEPE'; c='https://coinranking.com/coin/03WI8NQPF+pepe-pepe'; d='0.0000102'; e='  -6.97%  ' }
    @{ row=27; b='InternetComputer(DFINITY)'; c='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; d='8.80'; e='  -0.07%  ' }
    @{ row=28; b='Binance-PegBSC-USD'; c='https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; d='1.00'; e='  +0.07%  ' }
    @{ row=29; b='PancakeSwap'; c='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; d='1.81'; e='  -2.22%  ' }
    @{ row=30; b='EthereumClassic'; c='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; d='21.16'; e='  -3.12%  ' }
    @{ row=31; b='Aptos'; c='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; d='6.33'; e='  -0.92%  ' }
    @{ row=32; b='NEARProtocol'; c='https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; d='4.89'; e='  -2.11%  ' }
    @{ row=33; b='Fetch.AI'; c='https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; d='1.13'; e='  -1.18%  ' }
    @{ row=34; b='Monero'; c='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; d='156.61'; e='  +0.17%  ' }
    @{ row=35; b='ImmutableX'; c='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; d='1.34'; e='  -3.51%  ' }
    @{ row=36; b='Maker'; c='https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; d='2.696.36'; e='  -1.82%  ' }
    @{ row=37; b='Stacks'; c='https://coinranking.com/coin/mMPrMcB7+stacks-stx'; d='1.69'; e='  -1.02%  ' }
    @{ row=38; b='EnergySwap'; c='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; d='23.97'; e='  -4.62%  ' }
    @{ row=39; b='Filecoin'; c='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; d='4.07'; e='  -2.99%  ' }
    @{ row=40; b='Mantle'; c='https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; d='0.700'; e='  -2.60%  ' }
    @{ row=41; b='Hedera'; c='https://coinranking.com/coin/jad286TjB+hedera-hbar'; d='0.0621'; e='  -1.03%  ' }
    @{ row=42; b='RenderToken'; c='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; d='5.43'; e='  -6.40%  ' }
    @{ row=43; b='VeChain'; c='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; d='0.0257'; e='  -1.47%  ' }
    @{ row=44; b='Bittensor'; c='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; d='290.26'; e='  -2.65%  ' }
    @{ row=45; b='InjectiveProtocol'; c='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; d='21.35'; e='  -3.31%  ' }
    @{ row=46; b='FirstDigitalUSD'; c='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; d='1.00'; e='  +0.16%  ' }
    @{ row=47; b='Stellar'; c='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; d='0.0981'; e='  -1.36%  ' }
    @{ row=48; b='WhiteBITCoin'; c='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; d='10.48'; e='  +0.36%  ' }
    @{ row=49; b='dogwifhat'; c='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; d='1.89'; e='  -11.22%  ' }
    @{ row=50; b='Cosmos'; c='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; d='5.73'; e='  -1.21%  ' }
    @{ row=51; b='ONDO'; c='https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; d='0.880'; e='  -5.15%  ' }
)

foreach ($item in $data) {
    Set-TextValue $ws.Cells.Item($item.row, 2) $item.b
    Set-TextValue $ws.Cells.Item($item.row, 3) $item.c
    Set-TextValue $ws.Cells.Item($item.row, 4) $item.d
    Set-TextValue $ws.Cells.Item($item.row, 5) $item.e
}

Write-Output "Done updating $($data.Count) rows"
